$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.15552818775177
$ws.Range("B1").Value = 2.382137775421143
$ws.Range("D1").Value = 2.393799066543579
$ws.Range("E1").Value = 1.22511351108551
